# Multi browser implementation for chrome and firefox and accounting cash
# -------------------------------------------------------------------------
# 1. Normalise the "Currency" label/value pair on the ProductLoanInput
#    sheet (row 6): "Currency" -> "currency", "US Dollar " -> "US Dollar".
# 2. Drop the now-unused column C helper cells (C5, C6, C12).
# 3. Re-point the active sheet/selection: ProductLoanInput becomes the
#    active tab with A6:B6 selected; ProductLoanOutput keeps a B14
#    selection but is no longer the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoanInput
$ws2 = $wb.Worksheets.Item(2)   # ProductLoanOutput

# --- content edits -------------------------------------------------------
$ws1.Range("A6").Value = "currency"
$ws1.Range("B6").Value = "US Dollar"

# Remove the stray column C cells that accompanied the old layout.
$ws1.Range("C5").Clear()
$ws1.Range("C6").Clear()
$ws1.Range("C12").Clear()

# --- view/selection edits -------------------------------------------------
# Set the output sheet's selection first (it must not become the active tab).
$ws2.Range("B14").Select()

# Activate the input sheet and select A6:B6, clearing the old
# topLeftCell="A22" / F10 selection in the process.
$ws1.Activate()
$ws1.Range("A6:B6").Select()
